$wb = $excel.ActiveWorkbook

# Sheet ALC, row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 781.25
$ws.Range("I2").Value = 700
$ws.Range("J2").Value = 830
$ws.Range("K2").Value = 700
$ws.Range("L2").Value = 830
$ws.Range("M2").Value = -587
$ws.Range("N2").Value = -1056

# Sheet ALC, row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 312.16666
$ws.Range("J12").Value = 199.66667
$ws.Range("L12").Value = 199.66667
$ws.Range("N12").Value = -539.6666700000001

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 74999
$ws.Range("J18").Value = 74999
$ws.Range("L18").Value = 74999
$ws.Range("N18").Value = -75567

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1254.3928
$ws.Range("I33").Value = 1052.174
$ws.Range("K33").Value = 1052.174
$ws.Range("M33").Value = -823.174

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 70331.336
$ws.Range("J51").Value = 102747
$ws.Range("L51").Value = 102747
$ws.Range("N51").Value = -103715

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4579.476
$ws.Range("I137").Value = 1605.5714
$ws.Range("J137").Value = 10527.286
$ws.Range("K137").Value = 4816.7142
$ws.Range("L137").Value = 31581.858
$ws.Range("M137").Value = -2266.7142
$ws.Range("N137").Value = -36681.858

# Sheet ARM, row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 486.375
$ws.Range("I6").Value = 482.5
$ws.Range("K6").Value = 482.5
$ws.Range("M6").Value = -309.5

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 29479078
$ws.Range("I61").Value = 41673320
$ws.Range("K61").Value = 41673320
$ws.Range("M61").Value = -41673108

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3055.375
$ws.Range("I63").Value = 3055.375
$ws.Range("K63").Value = 3055.375
$ws.Range("M63").Value = -2369.375

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3055.375
$ws.Range("I66").Value = 3055.375
$ws.Range("K66").Value = 15276.875
$ws.Range("M66").Value = -11844.875

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7583644
$ws.Range("I74").Value = 10871541
$ws.Range("J74").Value = 21481.6
$ws.Range("K74").Value = 10871541
$ws.Range("L74").Value = 21481.6
$ws.Range("M74").Value = -10870667
$ws.Range("N74").Value = -23229.6

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7583644
$ws.Range("I77").Value = 10871541
$ws.Range("J77").Value = 21481.6
$ws.Range("K77").Value = 54357705
$ws.Range("L77").Value = 107408
$ws.Range("M77").Value = -54353337
$ws.Range("N77").Value = -116144

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1420.8182
$ws.Range("J122").Value = 1749.5
$ws.Range("L122").Value = 5248.5
$ws.Range("N122").Value = -10148.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6059.2173
$ws.Range("I132").Value = 4138.264
$ws.Range("K132").Value = 12414.792
$ws.Range("M132").Value = -9884.792000000001

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 29479078
$ws.Range("I136").Value = 41673320
$ws.Range("K136").Value = 125019960
$ws.Range("M136").Value = -125017410

# Sheet BSM, row 96
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 36666.2
$ws.Range("I96").Value = 14639.8
$ws.Range("J96").Value = 58692.6
$ws.Range("K96").Value = 14639.8
$ws.Range("L96").Value = 58692.6
$ws.Range("M96").Value = -11893.8
$ws.Range("N96").Value = -64184.6

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2062.2727
$ws.Range("I105").Value = 1627.9166
$ws.Range("J105").Value = 2583.5
$ws.Range("K105").Value = 1627.9166
$ws.Range("L105").Value = 2583.5
$ws.Range("M105").Value = 119.0834
$ws.Range("N105").Value = -6077.5

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33824.027
$ws.Range("I134").Value = 1778.85
$ws.Range("K134").Value = 5336.549999999999
$ws.Range("M134").Value = -2801.549999999999

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 557.5
$ws.Range("I22").Value = 557.5
$ws.Range("K22").Value = 557.5
$ws.Range("M22").Value = -207.5

# Sheet CRP, row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 49666.5
$ws.Range("J51").Value = 67999.8
$ws.Range("L51").Value = 67999.8
$ws.Range("N51").Value = -69471.8

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 73530.25
$ws.Range("I60").Value = 75121.5
$ws.Range("J60").Value = 72999.836
$ws.Range("K60").Value = 75121.5
$ws.Range("L60").Value = 72999.836
$ws.Range("M60").Value = -74610.5
$ws.Range("N60").Value = -74021.836

# Sheet CRP, row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 49666.5
$ws.Range("J61").Value = 67999.8
$ws.Range("L61").Value = 67999.8
$ws.Range("N61").Value = -68695.8

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5559716
$ws.Range("I4").Value = 15893504
$ws.Range("K4").Value = 47680512
$ws.Range("M4").Value = -47680400

# Sheet CUL, row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 666666
$ws.Range("J9").Value = 666666
$ws.Range("L9").Value = 1999998
$ws.Range("N9").Value = -2000446

# Sheet CUL, row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""

# Sheet CUL, row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""

# Sheet CUL, row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 82999.336
$ws.Range("J37").Value = 82999.336
$ws.Range("L37").Value = 248998.008
$ws.Range("N37").Value = -249222.008

# Sheet CUL, row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 422
$ws.Range("I40").Value = 558
$ws.Range("J40").Value = 150
$ws.Range("K40").Value = 2232
$ws.Range("L40").Value = 600
$ws.Range("M40").Value = -2163
$ws.Range("N40").Value = -738

# Sheet CUL, row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 702.6923
$ws.Range("J109").Value = 2153
$ws.Range("L109").Value = 6459
$ws.Range("N109").Value = -8539

# Sheet CUL, row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 947.9
$ws.Range("J117").Value = 620
$ws.Range("L117").Value = 1860
$ws.Range("N117").Value = -8744

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8578.308000000001
$ws.Range("J131").Value = 5658.8
$ws.Range("L131").Value = 16976.4
$ws.Range("N131").Value = -27056.4

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2674.95
$ws.Range("I139").Value = 2833
$ws.Range("J139").Value = 2647.0588
$ws.Range("K139").Value = 8499
$ws.Range("L139").Value = 7941.176399999999
$ws.Range("M139").Value = -3359
$ws.Range("N139").Value = -18221.1764

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 276560.62
$ws.Range("I140").Value = 432292.56
$ws.Range("J140").Value = 4029.75
$ws.Range("K140").Value = 1296877.68
$ws.Range("L140").Value = 12089.25
$ws.Range("M140").Value = -1291697.68
$ws.Range("N140").Value = -22449.25

# Sheet GSM, row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 73500
$ws.Range("J63").Value = 73500
$ws.Range("L63").Value = 73500
$ws.Range("N63").Value = -74872

# Sheet GSM, row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 73500
$ws.Range("J66").Value = 73500
$ws.Range("L66").Value = 220500
$ws.Range("N66").Value = -227364

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1644.8
$ws.Range("J80").Value = 1499.5
$ws.Range("L80").Value = 1499.5
$ws.Range("N80").Value = -3495.5

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1644.8
$ws.Range("J83").Value = 1499.5
$ws.Range("L83").Value = 7497.5
$ws.Range("N83").Value = -17481.5

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8181
$ws.Range("I102").Value = 4817.2
$ws.Range("K102").Value = 4817.2
$ws.Range("M102").Value = -3195.2

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3616.238
$ws.Range("I113").Value = 3566.3333
$ws.Range("J113").Value = 3741
$ws.Range("K113").Value = 3566.3333
$ws.Range("L113").Value = 3741
$ws.Range("M113").Value = -1396.3333
$ws.Range("N113").Value = -8081

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2340.6667
$ws.Range("I122").Value = 2340.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7022.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4572.000100000001
$ws.Range("N122").Value = ""

# Sheet LTW, row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828

# Sheet LTW, row 36
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 85853.5
$ws.Range("J36").Value = 85853.5
$ws.Range("L36").Value = 85853.5
$ws.Range("N36").Value = -86977.5

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3142.7856
$ws.Range("I40").Value = 1444.1111
$ws.Range("K40").Value = 1444.1111
$ws.Range("M40").Value = -1308.1111

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 52903.293
$ws.Range("I136").Value = 5483.3335
$ws.Range("J136").Value = 195163.17
$ws.Range("K136").Value = 16450.0005
$ws.Range("L136").Value = 585489.51
$ws.Range("M136").Value = -13900.0005
$ws.Range("N136").Value = -590589.51

# Sheet WVR, row 53
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 27617.4
$ws.Range("J53").Value = 27617.4
$ws.Range("L53").Value = 27617.4
$ws.Range("N53").Value = -28831.4

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17397.295
$ws.Range("I136").Value = 1156.8
$ws.Range("K136").Value = 3470.4
$ws.Range("M136").Value = -920.3999999999996
